# feat: add 2022-Q3 data
#
# 1. "总计" (summary) sheet: insert a new top data row for 2022-Q3 and
#    renumber the existing rows' running index (column A).
# 2. Duplicate the current "2022-Q2" sheet, drop the duplicate right
#    before the original, rename the duplicate "2022-Q3" and give it the
#    new quarter's figures. The original "2022-Q2" sheet (and the
#    "2021-Q2" / "2021-Q1" / "2020-Q4" sheets after it) keep their old
#    data untouched - they simply sit one tab further to the right now.

$wb = $excel.ActiveWorkbook

# --- 1. Update the "总计" summary sheet -----------------------------------
$summary = $wb.Worksheets.Item("总计")

# Make room for the new first data row (row 2); existing rows 2-5 shift
# down to rows 3-6.
$summary.Rows.Item(2).Insert()

# Seed A2 with row 3's formatting (border + centered index-column style)
# before overwriting the values, since Insert() leaves A2 unstyled.
$summary.Range("A3").Copy($summary.Range("A2"))

# New row 2: 2022-Q3
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.2
$summary.Range("B2:D2").Style = "Normal"

# Column A is a plain running index (0,1,2,...) - renumber rows 3-6 which
# used to be rows 2-5 (values 0,1,2,3) and are now one slot later
# (values 1,2,3,4).
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# --- 2. Add the "2022-Q3" fund sheet --------------------------------------
$oldQ2 = $wb.Worksheets.Item("2022-Q2")

# Duplicate "2022-Q2" and drop the copy immediately before the original -
# this gives us a ready-made sheet with identical layout/styles. (After
# Copy(), the duplicate gets auto-named "2022-Q2 (2)"; fetch sheets by
# name afterwards rather than trusting stale object references.)
$oldQ2.Copy($oldQ2)
$newQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$newQ3.Name = "2022-Q3"

# Overwrite the copied (old) figures with the new 2022-Q3 figures. These
# columns are stored as text in this workbook, so force text the same way
# Excel does (leading apostrophe) and then drop the resulting quote-prefix
# style so the cell ends up with no explicit style, matching the rest of
# the sheet.
$newQ3.Range("D2").Value = "'17.28"
$newQ3.Range("E2").Value = "'29.21"
$newQ3.Range("F2").Value = "'1.13"
$newQ3.Range("G2").Value = "'0.1953"
$newQ3.Range("D2:G2").Style = "Normal"

# Ranking column is a real number.
$newQ3.Range("H2").Value = 8

# Restore the originally-selected tab ("2020-Q4"), which the Copy() call
# above bumped.
$wb.Worksheets.Item("2020-Q4").Activate()
